$wb = $excel.ActiveWorkbook

# hunk 0: ALC row 12
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H12").Value = 11111981
$ws.Range("I12").Value = 25000084
$ws.Range("J12").Value = 1498.4
$ws.Range("K12").Value = 25000084
$ws.Range("L12").Value = 1498.4
$ws.Range("M12").Value = -24999914
$ws.Range("N12").Value = -1838.4

# hunk 1: ALC row 33
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H33").Value = 685.9091
$ws.Range("I33").Value = 788
$ws.Range("K33").Value = 788
$ws.Range("M33").Value = -559

# hunk 2: ALC row 99
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H99").Value = 203.5
$ws.Range("I99").Value = 188
$ws.Range("K99").Value = 564
$ws.Range("M99").Value = 934

# hunk 3: ALC row 101
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H101").Value = 3493.5386
$ws.Range("I101").Value = 602.3333
$ws.Range("K101").Value = 1806.9999
$ws.Range("M101").Value = -184.9999

# hunk 4: ALC row 103
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H103").Value = 393.5
$ws.Range("I103").Value = 299.8
$ws.Range("J103").Value = 549.6667
$ws.Range("K103").Value = 899.4000000000001
$ws.Range("L103").Value = 1649.0001
$ws.Range("M103").Value = -313.4000000000001
$ws.Range("N103").Value = -2821.0001

# hunk 5: ALC row 116
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H116").Value = 4763.091
$ws.Range("I116").Value = 4099
$ws.Range("K116").Value = 4099
$ws.Range("M116").Value = -657

# hunk 6: ALC row 118
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H118").Value = 197.5
$ws.Range("J118").Value = 0
$ws.Range("L118").Value = 0
$ws.Range("N118").ClearContents()

# hunk 7: ALC row 132
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H132").Value = 1768.619
$ws.Range("I132").Value = 1549.5264
$ws.Range("K132").Value = 4648.5792
$ws.Range("M132").Value = -2118.5792

# hunk 8: ALC row 135
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H135").Value = 561.72
$ws.Range("I135").Value = 480.33334
$ws.Range("K135").Value = 4323.00006
$ws.Range("M135").Value = -1788.00006

# hunk 9: ALC row 136
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H136").Value = 135311.81
$ws.Range("J136").Value = 135311.81
$ws.Range("L136").Value = 135311.81
$ws.Range("N136").Value = -145511.81

# hunk 10: ALC row 138
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H138").Value = 1991.6604
$ws.Range("I138").Value = 1294.5
$ws.Range("J138").Value = 2414.182
$ws.Range("K138").Value = 3883.5
$ws.Range("L138").Value = 7242.545999999999
$ws.Range("M138").Value = 1256.5
$ws.Range("N138").Value = -17522.546

# hunk 11: ARM row 32
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 1425.9275
$ws.Range("I32").Value = 1454.4395
$ws.Range("K32").Value = 1454.4395
$ws.Range("M32").Value = -1167.4395

# hunk 12: ARM row 45
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H45").Value = 6440.8335
$ws.Range("J45").Value = 7937.5
$ws.Range("L45").Value = 7937.5
$ws.Range("N45").Value = -8691.5

# hunk 13: ARM row 74
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H74").Value = 2777.04
$ws.Range("I74").Value = 2008.1892
$ws.Range("J74").Value = 4965.3076
$ws.Range("K74").Value = 2008.1892
$ws.Range("L74").Value = 4965.3076
$ws.Range("M74").Value = -1134.1892
$ws.Range("N74").Value = -6713.3076

# hunk 14: ARM row 77
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H77").Value = 2777.04
$ws.Range("I77").Value = 2008.1892
$ws.Range("J77").Value = 4965.3076
$ws.Range("K77").Value = 10040.946
$ws.Range("L77").Value = 24826.538
$ws.Range("M77").Value = -5672.946
$ws.Range("N77").Value = -33562.538

# hunk 15: ARM row 102
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H102").Value = 4832.4346
$ws.Range("I102").Value = 3819.7222
$ws.Range("J102").Value = 8478.200000000001
$ws.Range("K102").Value = 3819.7222
$ws.Range("L102").Value = 8478.200000000001
$ws.Range("M102").Value = -2197.7222
$ws.Range("N102").Value = -11722.2

# hunk 16: ARM row 106
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H106").Value = 0
$ws.Range("J106").Value = 0
$ws.Range("L106").Value = 0
$ws.Range("N106").ClearContents()

# hunk 17: BSM row 20
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 3796.0908
$ws.Range("I20").Value = 3361.3333
$ws.Range("J20").Value = 4097.077
$ws.Range("K20").Value = 3361.3333
$ws.Range("L20").Value = 4097.077
$ws.Range("M20").Value = -3114.3333
$ws.Range("N20").Value = -4591.077

# hunk 18: CRP row 58
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H58").Value = 11064.692
$ws.Range("I58").Value = 8472.5
$ws.Range("J58").Value = 12216.777
$ws.Range("K58").Value = 8472.5
$ws.Range("L58").Value = 12216.777
$ws.Range("M58").Value = -8269.5
$ws.Range("N58").Value = -12622.777

# hunk 19: CRP row 86
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H86").Value = 92565.5
$ws.Range("I86").Value = 10747
$ws.Range("K86").Value = 10747
$ws.Range("M86").Value = -9624

# hunk 20: CRP row 89
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H89").Value = 92565.5
$ws.Range("I89").Value = 10747
$ws.Range("K89").Value = 53735
$ws.Range("M89").Value = -48119

# hunk 21: CRP row 99
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H99").Value = 3525.5
$ws.Range("J99").Value = 2988.5715
$ws.Range("L99").Value = 2988.5715
$ws.Range("N99").Value = -5984.5715

# hunk 22: CRP row 126
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H126").Value = 3525.5
$ws.Range("J126").Value = 2988.5715
$ws.Range("L126").Value = 8965.7145
$ws.Range("N126").Value = -13905.7145

# hunk 23: CRP row 132
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H132").Value = 7500
$ws.Range("I132").Value = 7500
$ws.Range("J132").Value = 0
$ws.Range("K132").Value = 22500
$ws.Range("L132").Value = 0
$ws.Range("M132").Value = -19970
$ws.Range("N132").ClearContents()

# hunk 24: CRP row 136
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H136").Value = 11064.692
$ws.Range("I136").Value = 8472.5
$ws.Range("J136").Value = 12216.777
$ws.Range("K136").Value = 25417.5
$ws.Range("L136").Value = 36650.331
$ws.Range("M136").Value = -22867.5
$ws.Range("N136").Value = -41750.331

# hunk 25: CUL row 4
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H4").Value = 96033210
$ws.Range("I4").Value = 96033210
$ws.Range("K4").Value = 288099630
$ws.Range("M4").Value = -288099518

# hunk 26: CUL row 11
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H11").Value = 20000400
$ws.Range("I11").Value = 33333400
$ws.Range("J11").Value = 901
$ws.Range("K11").Value = 100000200
$ws.Range("L11").Value = 2703
$ws.Range("M11").Value = -100000060
$ws.Range("N11").Value = -2983

# hunk 27: CUL row 114
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H114").Value = 2500
$ws.Range("J114").Value = 2500
$ws.Range("L114").Value = 7500
$ws.Range("N114").Value = -14008

# hunk 28: CUL row 121
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H121").Value = 6733928.5
$ws.Range("I121").Value = 490.7
$ws.Range("K121").Value = 1472.1
$ws.Range("M121").Value = -162.0999999999999

# hunk 29: CUL row 129
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H129").Value = 19133842
$ws.Range("I129").Value = 33434536
$ws.Range("J129").Value = 1257974.8
$ws.Range("K129").Value = 100303608
$ws.Range("L129").Value = 3773924.4
$ws.Range("M129").Value = -100298608
$ws.Range("N129").Value = -3783924.4

# hunk 30: CUL row 134
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H134").Value = 1499.6666
$ws.Range("I134").Value = 999.5
$ws.Range("K134").Value = 2998.5
$ws.Range("M134").Value = 2071.5

# hunk 31: CUL row 140
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H140").Value = 435243.6
$ws.Range("I140").Value = 1056.1904
$ws.Range("K140").Value = 3168.5712
$ws.Range("M140").Value = 2011.4288

# hunk 32: GSM row 2
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H2").Value = 9221.182000000001
$ws.Range("J2").Value = 50012
$ws.Range("L2").Value = 50012
$ws.Range("N2").Value = -50238

# hunk 33: GSM row 36
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H36").Value = 5017
$ws.Range("I36").Value = 5017
$ws.Range("K36").Value = 5017
$ws.Range("M36").Value = -4532

# hunk 34: GSM row 46
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H46").Value = 46490.715
$ws.Range("J46").Value = 51608.75
$ws.Range("L46").Value = 51608.75
$ws.Range("N46").Value = -51920.75

# hunk 35: GSM row 122
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H122").Value = 4232.2666
$ws.Range("I122").Value = 5121.625
$ws.Range("J122").Value = 3215.8572
$ws.Range("K122").Value = 15364.875
$ws.Range("L122").Value = 9647.571599999999
$ws.Range("M122").Value = -12914.875
$ws.Range("N122").Value = -14547.5716

# hunk 36: GSM row 130
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H130").Value = 49999
$ws.Range("I130").Value = 0
$ws.Range("J130").Value = 49999
$ws.Range("K130").Value = 0
$ws.Range("L130").Value = 49999
$ws.Range("M130").ClearContents()
$ws.Range("N130").Value = -60039

# hunk 37: GSM row 132
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H132").Value = 2084.7144
$ws.Range("I132").Value = 2399.75
$ws.Range("K132").Value = 7199.25
$ws.Range("M132").Value = -4669.25

# hunk 38: LTW row 40
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 1958.3334
$ws.Range("I40").Value = 1958.3334
$ws.Range("K40").Value = 1958.3334
$ws.Range("M40").Value = -1822.3334

# hunk 39: LTW row 55
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H55").Value = 689.875
$ws.Range("J55").Value = 558.1667
$ws.Range("L55").Value = 558.1667
$ws.Range("N55").Value = -904.1667

# hunk 40: LTW row 122
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H122").Value = 4459.6
$ws.Range("I122").Value = 4449.5
$ws.Range("J122").Value = 4466.3335
$ws.Range("K122").Value = 13348.5
$ws.Range("L122").Value = 13399.0005
$ws.Range("M122").Value = -10898.5
$ws.Range("N122").Value = -18299.0005

# hunk 41: LTW row 132
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H132").Value = 15275.632
$ws.Range("I132").Value = 17351.785
$ws.Range("J132").Value = 9462.4
$ws.Range("K132").Value = 52055.355
$ws.Range("L132").Value = 28387.2
$ws.Range("M132").Value = -49525.355
$ws.Range("N132").Value = -33447.2

# hunk 42: WVR row 122
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H122").Value = 2523.652
$ws.Range("I122").Value = 1611.3572
$ws.Range("J122").Value = 3942.7778
$ws.Range("K122").Value = 4834.071599999999
$ws.Range("L122").Value = 11828.3334
$ws.Range("M122").Value = -2384.071599999999
$ws.Range("N122").Value = -16728.3334

# hunk 43: WVR row 132
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 6357.7393
$ws.Range("I132").Value = 5601.6113
$ws.Range("K132").Value = 16804.8339
$ws.Range("M132").Value = -14274.8339

# hunk 44: WVR row 136
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H136").Value = 4025.4285
$ws.Range("I136").Value = 2870.2632
$ws.Range("K136").Value = 8610.7896
$ws.Range("M136").Value = -6060.7896
